$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data for rows 2-31 (x, y, z), after inserting 4 rows before the
# old row 2 and appending 6 rows after the old row 21.
$finalData = @(
    @(-0.2807844198983286, -0.8123789481047923, -0.09823161813205905),
    @(-0.1417398627462046, -0.1100757933896112, 2.834019057966533),
    @(-1.13705500035451, 2.739093329500554, -2.630493768330314),
    @(2.9049899660308, -3.484631538391146, -2.321273335095107),
    @(-3.415270102435171, 0.5314338042818449, 4.029043727907666),
    @(-9.314493606830482, -0.2480275487077682, 3.206530965607738),
    @(2.801971287562949, -2.823644173042532, 0.2525709908584024),
    @(4.358451415752542, -4.83039114393037, 0.2956774553348279),
    @(7.929245866578198, 4.034436612293629, 0.2522587591204151),
    @(-0.834473840121543, 4.45597813047212, 2.253706974202184),
    @(-5.379544188236389, 2.00216909318133, 4.38414403076831),
    @(-1.495801197043751, 11.70056653022773, 0.6984025085794245),
    @(2.489601142961419, 1.782228860361755, 3.894388745057194),
    @(7.727643826912207, 6.834787146798529, 1.897908695812871),
    @(0.1587105783922702, 5.019907260763223, 1.782309127265012),
    @(-5.217137336730959, 3.712993112103693, 5.708451303942454),
    @(1.931081188136128, 1.861609384931368, -5.062639532418077),
    @(3.396664668773761, -5.718913538702679, -1.399837789864275),
    @(6.86892145666587, -1.869332508794124, 1.079519285724066),
    @(1.767213952952357, 1.710883551630483, 1.963973425585661),
    @(-3.255076243959633, 0.04516811000888521, 5.548670966049739),
    @(-0.7715866319064268, 8.701057298430134, -0.006967774752950184),
    @(1.181180877932187, -0.4844702030050883, 1.036063860202677),
    @(5.493912943478295, -2.113819895119489, 1.030801805956604),
    @(2.463171288884921, 2.284831554725252, 0.6873917538544219),
    @(-3.070133529860395, 3.359637453638269, 3.026170623713525),
    @(-2.565556277488823, 10.97180558895251, 1.010093347779601),
    @(-0.5303593524571014, 2.632135161037951, 0.3993937804781349),
    @(4.747484356164973, -6.67747688293454, -1.334976077079803),
    @(2.688300560260526, -1.088815121815081, -1.408782399933858)
)

for ($i = 0; $i -lt $finalData.Length; $i++) {
    $r = 2 + $i
    $row = $finalData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
